# Apply the "New simulation files for schemes report" edit.
# - Row 2 (header) HKL labels C2:J2 are reordered.
# - Existing data rows 3-19 column B get relabeled (shared string table
#   reorder/replacement means the scheme names shown change even though
#   the row index/count values stay the same).
# - 10 new data rows (20-29) are appended for the new schemes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the reordered HKL column headers ---
$ws.Range("C2").Value = "[3, 2, 1]"
$ws.Range("D2").Value = "[3, 1, 0]"
$ws.Range("E2").Value = "[2, 2, 2]"
$ws.Range("F2").Value = "[1, 1, 0]"
$ws.Range("G2").Value = "[2, 0, 0]"
$ws.Range("H2").Value = "[2, 2, 0]"
$ws.Range("I2").Value = "[4, 0, 0]"
$ws.Range("J2").Value = "[2, 1, 1]"

# --- Rows 3-19: relabel column B (scheme names) ---
$ws.Range("B3").Value = "Spiral5"
$ws.Range("B4").Value = "RotRing OmegaMax-90"
$ws.Range("B5").Value = "Equal Angle"
$ws.Range("B6").Value = "Tilt Rotate"
$ws.Range("B7").Value = "CLR"
$ws.Range("B8").Value = "Rizzie Hex"
$ws.Range("B9").Value = "Thomas Hex"
$ws.Range("B10").Value = "Tilt Rotate_Partial"
$ws.Range("B11").Value = "RotRing OmegaMax-60"
$ws.Range("B12").Value = "Equal Angle_Partial"
$ws.Range("B13").Value = "Rizzie Hex_Partial"
$ws.Range("B14").Value = "ND Single"
$ws.Range("B15").Value = "RD Single"
$ws.Range("B16").Value = "TD Single"
$ws.Range("B17").Value = "Morris Single"
$ws.Range("B18").Value = "Ring Perpendicular to ND"
$ws.Range("B19").Value = "Ring Perpendicular to RD"

# --- Rows 20-29: new data rows, formatted like the existing block ---
$ws.Range("A19:T19").Copy()
$ws.Range("A20:T29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Ring Perpendicular to TD"
$ws.Range("C20:T20").Value = 1
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "OffsetFTD"
$ws.Range("C21:T21").Value = 1
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "OffsetATD"
$ws.Range("C22:T22").Value = 1
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "OffsetF45"
$ws.Range("C23:T23").Value = 1
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "OffsetA45"
$ws.Range("C24:T24").Value = 1
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "OffsetFRD"
$ws.Range("C25:T25").Value = 1
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "OffsetARD"
$ws.Range("C26:T26").Value = 1
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "Gaussian Quadrature"
$ws.Range("C27:T27").Value = 1
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "Michael-CCHex"
$ws.Range("C28:T28").Value = 1
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "Michael-SNHex"
$ws.Range("C29:T29").Value = 1

$ws.Range("A1").Select()
